$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C6: port value becomes combined "443,80" (stored as text)
$ws.Range("C6").Value = "443,80"

# Update F6: tag value becomes combined "machine,servers"
$ws.Range("F6").Value = "machine,servers"

# Update the active selection to C7
$ws.Range("C7").Select()
